{"js": "// Each entry is [originalEquation, newEquation]; every value in the\n// document's multiplication-practice table is unique, so a plain text\n// search+replace for each pair is unambiguous.\nconst replacements = [\n  [\"98\u00d734=3332\", \"34\u00d778=2652\"],\n  [\"59\u00d765=3835\", \"50\u00d772=3600\"],\n  [\"28\u00d755=1540\", \"33\u00d741=1353\"],\n  [\"26\u00d778=2028\", \"62\u00d783=5146\"],\n  [\"90\u00d748=4320\", \"84\u00d789=7476\"],\n  [\"21\u00d754=1134\", \"30\u00d779=2370\"],\n  [\"18\u00d743=774\", \"85\u00d789=7565\"],\n  [\"72\u00d734=2448\", \"81\u00d775=6075\"],\n  [\"31\u00d711=341\", \"21\u00d730=630\"],\n  [\"40\u00d757=2280\", \"44\u00d719=836\"],\n  [\"48\u00d798=4704\", \"87\u00d793=8091\"],\n  [\"70\u00d792=6440\", \"15\u00d722=330\"],\n  [\"64\u00d763=4032\", \"57\u00d714=798\"],\n  [\"21\u00d792=1932\", \"16\u00d728=448\"],\n  [\"45\u00d780=3600\", \"50\u00d725=1250\"],\n  [\"67\u00d759=3953\", \"57\u00d772=4104\"],\n  [\"35\u00d779=2765\", \"86\u00d722=1892\"],\n  [\"71\u00d790=6390\", \"82\u00d795=7790\"],\n  [\"71\u00d760=4260\", \"25\u00d784=2100\"],\n  [\"92\u00d754=4968\", \"32\u00d736=1152\"],\n  [\"97\u00d723=2231\", \"33\u00d799=3267\"],\n  [\"79\u00d798=7742\", \"98\u00d765=6370\"],\n  [\"23\u00d785=1955\", \"22\u00d719=418\"],\n  [\"36\u00d753=1908\", \"88\u00d756=4928\"],\n  [\"79\u00d713=1027\", \"91\u00d797=8827\"],\n];\n\nconst body = context.document.body;\nconst searchResults = [];\n\n// Queue up a search for each original equation.\nfor (const [find] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  searchResults.push(results);\n}\n\nawait context.sync();\n\n// Replace each found occurrence with its new equation.\nlet replacedCount = 0;\nfor (let i = 0; i < replacements.length; i++) {\n  const [, replaceWith] = replacements[i];\n  const results = searchResults[i];\n  for (const item of results.items) {\n    item.insertText(replaceWith, Word.InsertLocation.replace);\n    replacedCount++;\n  }\n}\n\nawait context.sync();\n\nreturn `replaced ${replacedCount} of ${replacements.length} patterns`;\n", "ps1": "# Each entry is (originalEquation, newEquation); every value in the\n# document's multiplication-practice table is unique, so a plain text\n# Find/Replace for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"98\u00d734=3332\", \"34\u00d778=2652\"),\n    @(\"59\u00d765=3835\", \"50\u00d772=3600\"),\n    @(\"28\u00d755=1540\", \"33\u00d741=1353\"),\n    @(\"26\u00d778=2028\", \"62\u00d783=5146\"),\n    @(\"90\u00d748=4320\", \"84\u00d789=7476\"),\n    @(\"21\u00d754=1134\", \"30\u00d779=2370\"),\n    @(\"18\u00d743=774\", \"85\u00d789=7565\"),\n    @(\"72\u00d734=2448\", \"81\u00d775=6075\"),\n    @(\"31\u00d711=341\", \"21\u00d730=630\"),\n    @(\"40\u00d757=2280\", \"44\u00d719=836\"),\n    @(\"48\u00d798=4704\", \"87\u00d793=8091\"),\n    @(\"70\u00d792=6440\", \"15\u00d722=330\"),\n    @(\"64\u00d763=4032\", \"57\u00d714=798\"),\n    @(\"21\u00d792=1932\", \"16\u00d728=448\"),\n    @(\"45\u00d780=3600\", \"50\u00d725=1250\"),\n    @(\"67\u00d759=3953\", \"57\u00d772=4104\"),\n    @(\"35\u00d779=2765\", \"86\u00d722=1892\"),\n    @(\"71\u00d790=6390\", \"82\u00d795=7790\"),\n    @(\"71\u00d760=4260\", \"25\u00d784=2100\"),\n    @(\"92\u00d754=4968\", \"32\u00d736=1152\"),\n    @(\"97\u00d723=2231\", \"33\u00d799=3267\"),\n    @(\"79\u00d798=7742\", \"98\u00d765=6370\"),\n    @(\"23\u00d785=1955\", \"22\u00d719=418\"),\n    @(\"36\u00d753=1908\", \"88\u00d756=4928\"),\n    @(\"79\u00d713=1027\", \"91\u00d797=8827\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue=1, wdReplaceOne=2 (only one match exists per string anyway)\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
